$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("used")

# The two names that were "used" (picked off the top of Sheet1's pool,
# rows 2 and 3) get logged into the "used" sheet with their source file
# and the time they were used, then removed from the available pool.
$usedNames = @(
    @{ Id = "0rh82tbc"; File = "ChatGPT Image 2026年1月21日 16_55_57.png"; UsedAt = "2026-01-21 17:00:07" },
    @{ Id = "r88cbjqo"; File = "ChatGPT Image 2026年1月21日 16_58_50.png"; UsedAt = "2026-01-21 17:00:07" }
)

# Next empty row in the "used" log sheet (currently rows 1-36 are filled).
$nextRow = $ws2.Cells.Item($ws2.Rows.Count, 1).End(-4162).Row + 1

foreach ($entry in $usedNames) {
    $ws2.Cells.Item($nextRow, 1).Value = $entry.Id
    $ws2.Cells.Item($nextRow, 2).Value = $entry.File
    $ws2.Cells.Item($nextRow, 3).Value = $entry.UsedAt
    $nextRow = $nextRow + 1

    # Remove the consumed name from the front of the Sheet1 pool (row 2 is
    # always the next candidate once the header/first row stays put).
    $ws1.Rows.Item(2).Delete()
}
